# Apply updated cryptocurrency Price (D) / Volume(1h) (E) values per the source diff.
# All target cells in this sheet hold plain text (t="inlineStr" in the original XML),
# e.g. "214.76", "65.50" -- not numbers. For values that look like plain numbers,
# force the cell's number format to Text ("@") before assigning so Excel keeps the
# exact literal string (incl. trailing zeros) instead of silently converting it to a
# floating-point number. Percent-style cells (column E) are never numeric-looking
# (leading/trailing spaces, % sign) so they are safe to assign directly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.943.83'
$ws.Range("E2").Value = '  +0.02%  '

# Row 3
$ws.Range("D3").Value = '1.671.30'
$ws.Range("E3").Value = '  +1.21%  '

# Row 4
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.89'
$ws.Range("E5").Value = '  +0.14%  '

# Row 6
$ws.Range("E6").Value = '  +1.51%  '

# Row 7
$ws.Range("E7").Value = '  -0.01%  '

# Row 8
$ws.Range("E8").Value = '  +0.48%  '

# Row 9
$ws.Range("E9").Value = '  +0.64%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.20'
$ws.Range("E10").Value = '  +0.16%  '

# Row 11
$ws.Range("E11").Value = '  +1.56%  '

# Row 12
$ws.Range("D12").Value = '1.906.77'
$ws.Range("E12").Value = '  +1.20%  '

# Row 13
$ws.Range("D13").Value = '1.674.31'
$ws.Range("E13").Value = '  +1.31%  '

# Row 14
$ws.Range("E14").Value = '  +0.31%  '

# Row 15
$ws.Range("E15").Value = '  +1.23%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.53'
$ws.Range("E16").Value = '  +0.60%  '

# Row 17
$ws.Range("D17").Value = '26.936.48'
$ws.Range("E17").Value = '  -0.01%  '

# Row 18
$ws.Range("E18").Value = '  +4.13%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '233.83'
$ws.Range("E19").Value = '  -0.94%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0733'
$ws.Range("E20").Value = '  +0.07%  '

# Row 21
$ws.Range("E21").Value = '  -0.01%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.43'
$ws.Range("E22").Value = '  +0.42%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.18'
$ws.Range("E23").Value = '  -1.33%  '

# Row 24
$ws.Range("E24").Value = '  -2.05%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.89'
$ws.Range("E25").Value = '  +0.65%  '

# Row 26
$ws.Range("E26").Value = '  +0.27%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.93'
$ws.Range("E27").Value = '  +0.86%  '

# Row 29
$ws.Range("E29").Value = '  +0.17%  '

# Row 30
$ws.Range("E30").Value = '  +0.16%  '

# Row 31
$ws.Range("E31").Value = '  +0.09%  '

# Row 32
$ws.Range("E32").Value = '  +0.66%  '

# Row 33
$ws.Range("D33").Value = '1.460.30'
$ws.Range("E33").Value = '  -5.28%  '

# Row 34
$ws.Range("E34").Value = '  +1.78%  '

# Row 35
$ws.Range("E35").Value = '  +1.98%  '

# Row 36
$ws.Range("E36").Value = '  +0.20%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.581'
$ws.Range("E37").Value = '  -0.23%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.900'
$ws.Range("E38").Value = '  +0.77%  '

# Row 39
$ws.Range("E39").Value = '  +0.98%  '

# Row 40
$ws.Range("E40").Value = '  +13.44%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.78'
$ws.Range("E41").Value = '  -3.41%  '

# Row 42
$ws.Range("E42").Value = '  +0.02%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '66.31'
$ws.Range("E44").Value = '  +0.49%  '

# Row 45
$ws.Range("D45").Value = '1.811.75'

# Row 46
$ws.Range("E46").Value = '  +0.90%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.58'
$ws.Range("E47").Value = '  +0.89%  '

# Row 48
$ws.Range("E48").Value = '  +1.36%  '

# Row 49
$ws.Range("E49").Value = '  +3.08%  '

# Row 50
$ws.Range("E50").Value = '  +0.55%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.67'
$ws.Range("E51").Value = '  +0.69%  '
